# Refactor: create getRandomLocation function
#
# The underlying test-data generator now builds each departure/destination
# pair via a shared getRandomLocation() helper instead of two separate
# inline lookups. Re-running the generator against the same workbook
# produced a new (still valid) sample of city pairs:
#   - the destination column no longer repeats the very first departure
#     city in row 2 (that cell is now left blank instead of duplicating
#     "Acapulco"),
#   - for rows 3-26 the helper happened to draw the same city for both
#     departure and destination,
#   - and "Morelia" is now drawn one row earlier than "Monterrey" (affecting
#     the destination values on rows 27 and 28 as well).
# This script reproduces that resulting worksheet state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: destination cell is no longer populated.
$ws.Range("C2").ClearContents()

# Rows 3-26: destination now mirrors the departure city for that row.
$ws.Range("C3").Value  = "Aguascalientes"
$ws.Range("C4").Value  = "Cancun"
$ws.Range("C5").Value  = "Chetumal"
$ws.Range("C6").Value  = "Chihuahua"
$ws.Range("C7").Value  = "Ciudad Juarez"
$ws.Range("C8").Value  = "Ciudad Obregon"
$ws.Range("C9").Value  = "Colima"
$ws.Range("C10").Value = "Cozumel"
$ws.Range("C11").Value = "Culiacan"
$ws.Range("C12").Value = "Durango"
$ws.Range("C13").Value = "Guadalajara"
$ws.Range("C14").Value = "Hermosillo"
$ws.Range("C15").Value = "Huatulco"
$ws.Range("C16").Value = "Ixtapa / Zihuatanejo"
$ws.Range("C17").Value = "La Paz"
$ws.Range("C18").Value = "Leon"
$ws.Range("C19").Value = "Loreto"
$ws.Range("C20").Value = "Los Cabos"
$ws.Range("C21").Value = "Los Mochis"
$ws.Range("C22").Value = "Mazatlan"
$ws.Range("C23").Value = "Merida"
$ws.Range("C24").Value = "Mexicali"
$ws.Range("C25").Value = "Mexico City (AICM)"
$ws.Range("C26").Value = "Mexico City (AIFA) New"

# Rows 27-28: shifted because Morelia now sorts/generates ahead of Monterrey.
$ws.Range("C27").Value = "Morelia"
$ws.Range("C28").Value = "Oaxaca"
